$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 97, pushing current rows 97-106 down to 98-107.
$ws.Rows.Item(97).Insert()

# Populate the new row 97 with values (same pattern as surrounding rows for A,B,C,E,F,G,H,I,J).
$ws.Cells.Item(97, 1).Value = 1
$ws.Cells.Item(97, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(97, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(97, 4).Value = 44748
$ws.Cells.Item(97, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(97, 5).Value = 15
$ws.Cells.Item(97, 6).Value = "Fruta"
$ws.Cells.Item(97, 7).Value = 100106
$ws.Cells.Item(97, 8).Value = "Oleaginosos"
$ws.Cells.Item(97, 9).Value = 100106002
$ws.Cells.Item(97, 10).Value = "Palta"
$ws.Cells.Item(97, 11).Value = "Fuerte"
$ws.Cells.Item(97, 12).Value = "Tercera"
$ws.Cells.Item(97, 13).Value = 200
$ws.Cells.Item(97, 14).Value = 42000
$ws.Cells.Item(97, 15).Value = 45000
$ws.Cells.Item(97, 16).Value = 43500
$ws.Cells.Item(97, 17).Value = "$/caja 25 kilos"
$ws.Cells.Item(97, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(97, 19).Value = 1740
$ws.Cells.Item(97, 20).Value = 25
